$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D "Price" cells hold numeric-looking text (e.g. "308.83", "10.60",
# "26.917.26") that must stay plain text exactly as scraped, including
# trailing zeros and the thousands-dot formatting. Assigning such a string
# straight to .Value lets Excel auto-detect it as a number and normalize it
# (dropping trailing zeros, turning dotted big numbers into real numbers,
# etc.), so we briefly force the cell to Text format, write the value, then
# clear the formatting again so the cell style stays at its original default.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.917.26"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.815.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4675"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3687"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.51%  "

$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8702"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.18%  "

$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.802.03"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.373"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.515"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07061"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.61"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008708"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("E20").Value = "  -0.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.943.79"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.321"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.60"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.030.52"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.34%  "

$ws.Range("E25").Value = "  -0.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.45"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.168"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.33"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.333"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.26"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08931"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7674"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.162"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.504"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.926"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.28%  "

$ws.Range("E36").Value = "  +0.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.086"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01964"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05287"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.932"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.256"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5310"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.343"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.75%  "

$ws.Range("E44").Value = "  -0.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.434"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4919"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.42"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.55%  "

$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.670"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.80"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.48%  "

$ws.Range("E51").Value = "  -0.31%  "
